$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 202.25
$ws.Range("J6").Value = 20
$ws.Range("L6").Value = 60
$ws.Range("N6").Value = -284
$ws.Range("H8").Value = 13.666667
$ws.Range("I8").Value = 13.666667
$ws.Range("K8").Value = 41.000001
$ws.Range("M8").Value = 97.999999
$ws.Range("H11").Value = 113.875
$ws.Range("I11").Value = 113.875
$ws.Range("K11").Value = 113.875
$ws.Range("M11").Value = 26.125
$ws.Range("H17").Value = 409.4737
$ws.Range("I17").Value = 455
$ws.Range("J17").Value = 406.94446
$ws.Range("K17").Value = 1365
$ws.Range("L17").Value = 1220.83338
$ws.Range("M17").Value = -1197
$ws.Range("N17").Value = -1556.83338
$ws.Range("I27").Value = 10000
$ws.Range("K27").Value = 30000
$ws.Range("M27").Value = -29899
$ws.Range("H31").Value = 28099.8
$ws.Range("I31").Value = 28099.8
$ws.Range("K31").Value = 84299.39999999999
$ws.Range("M31").Value = -84069.39999999999
$ws.Range("H38").Value = 585.2308
$ws.Range("I38").Value = 139.77777
$ws.Range("J38").Value = 1587.5
$ws.Range("K38").Value = 419.33331
$ws.Range("L38").Value = 4762.5
$ws.Range("M38").Value = -47.33330999999998
$ws.Range("N38").Value = -5506.5
$ws.Range("H39").Value = 51.92857
$ws.Range("I39").Value = 32.333332
$ws.Range("J39").Value = 87.2
$ws.Range("K39").Value = 96.999996
$ws.Range("L39").Value = 261.6
$ws.Range("M39").Value = 199.000004
$ws.Range("N39").Value = -853.6
$ws.Range("H42").Value = 46.8
$ws.Range("I42").Value = 46
$ws.Range("K42").Value = 138
$ws.Range("M42").Value = 92
$ws.Range("H93").Value = 29900
$ws.Range("J93").Value = 29900
$ws.Range("L93").Value = 29900
$ws.Range("N93").Value = -34892
$ws.Range("H99").Value = 779
$ws.Range("J99").Value = 989
$ws.Range("L99").Value = 2967
$ws.Range("N99").Value = -5963
$ws.Range("H106").Value = 34500436
$ws.Range("I106").Value = 37053024
$ws.Range("K106").Value = 37053024
$ws.Range("M106").Value = -37052393
$ws.Range("H118").Value = 337.8
$ws.Range("I118").Value = 337.8
$ws.Range("K118").Value = 1013.4
$ws.Range("M118").Value = 643.5999999999999
$ws.Range("H127").Value = 5747.8335
$ws.Range("J127").Value = 1450
$ws.Range("L127").Value = 4350
$ws.Range("N127").Value = -14270
$ws.Range("H129").Value = 2321.1667
$ws.Range("J129").Value = 2815
$ws.Range("L129").Value = 8445
$ws.Range("N129").Value = -18445
$ws.Range("H132").Value = 2759.9473
$ws.Range("I132").Value = 1320
$ws.Range("K132").Value = 3960
$ws.Range("M132").Value = -1430
$ws.Range("H138").Value = 2176.3137
$ws.Range("I138").Value = 1975.7273
$ws.Range("K138").Value = 5927.1819
$ws.Range("M138").Value = -787.1818999999996
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 10495.25
$ws.Range("J44").Value = 10495.25
$ws.Range("L44").Value = 10495.25
$ws.Range("N44").Value = -11471.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2998.5
$ws.Range("I86").Value = 2998.5
$ws.Range("K86").Value = 2998.5
$ws.Range("M86").Value = -1875.5
$ws.Range("H89").Value = 2998.5
$ws.Range("I89").Value = 2998.5
$ws.Range("K89").Value = 14992.5
$ws.Range("M89").Value = -9376.5
$ws.Range("H94").Value = 1618.3636
$ws.Range("I94").Value = 509.66666
$ws.Range("J94").Value = 2385.923
$ws.Range("K94").Value = 509.66666
$ws.Range("L94").Value = 2385.923
$ws.Range("M94").Value = -58.66665999999998
$ws.Range("N94").Value = -3287.923
$ws.Range("H105").Value = 3972925
$ws.Range("I105").Value = 7580139
$ws.Range("K105").Value = 7580139
$ws.Range("M105").Value = -7578392

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 110186.8
$ws.Range("I22").Value = 111318.555
$ws.Range("K22").Value = 111318.555
$ws.Range("M22").Value = -110968.555
$ws.Range("H31").Value = 3229.8064
$ws.Range("I31").Value = 2362.8262
$ws.Range("K31").Value = 2362.8262
$ws.Range("M31").Value = -2067.8262
$ws.Range("H34").Value = 3229.8064
$ws.Range("I34").Value = 2362.8262
$ws.Range("K34").Value = 2362.8262
$ws.Range("M34").Value = -2160.8262
$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 1000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 1000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -1588
$ws.Range("H86").Value = 6639.1665
$ws.Range("I86").Value = 6639.1665
$ws.Range("K86").Value = 6639.1665
$ws.Range("M86").Value = -5516.1665
$ws.Range("H89").Value = 6639.1665
$ws.Range("I89").Value = 6639.1665
$ws.Range("K89").Value = 33195.8325
$ws.Range("M89").Value = -27579.8325
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 340.125
$ws.Range("I26").Value = 130
$ws.Range("J26").Value = 466.2
$ws.Range("K26").Value = 390
$ws.Range("L26").Value = 1398.6
$ws.Range("M26").Value = -102
$ws.Range("N26").Value = -1974.6
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H69").Value = 2199.5715
$ws.Range("J69").Value = 2149.5
$ws.Range("L69").Value = 6448.5
$ws.Range("N69").Value = -8070.5
$ws.Range("H72").Value = 2199.5715
$ws.Range("J72").Value = 2149.5
$ws.Range("L72").Value = 19345.5
$ws.Range("N72").Value = -27457.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 469.83334
$ws.Range("I2").Value = 104.875
$ws.Range("J2").Value = 1199.75
$ws.Range("K2").Value = 104.875
$ws.Range("L2").Value = 1199.75
$ws.Range("M2").Value = 8.125
$ws.Range("N2").Value = -1425.75
$ws.Range("H20").Value = 34629.332
$ws.Range("J20").Value = 34629.332
$ws.Range("L20").Value = 34629.332
$ws.Range("N20").Value = -35119.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12751
$ws.Range("I7").Value = 12751
$ws.Range("K7").Value = 12751
$ws.Range("M7").Value = -12639
$ws.Range("H13").Value = 100000000
$ws.Range("I13").Value = 100000000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 100000000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -99999860
$ws.Range("N13").ClearContents()
$ws.Range("H46").Value = 2895.5715
$ws.Range("I46").Value = 2223
$ws.Range("J46").Value = 3400
$ws.Range("K46").Value = 2223
$ws.Range("L46").Value = 3400
$ws.Range("M46").Value = -2035
$ws.Range("N46").Value = -3776
$ws.Range("H68").Value = 7779.6
$ws.Range("I68").Value = 8473.75
$ws.Range("K68").Value = 8473.75
$ws.Range("M68").Value = -7724.75
$ws.Range("H71").Value = 7779.6
$ws.Range("I71").Value = 8473.75
$ws.Range("K71").Value = 42368.75
$ws.Range("M71").Value = -38624.75
$ws.Range("H126").Value = 12751
$ws.Range("I126").Value = 12751
$ws.Range("K126").Value = 38253
$ws.Range("M126").Value = -35783

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3875.3333
$ws.Range("I62").Value = 3104
$ws.Range("K62").Value = 3104
$ws.Range("M62").Value = -2480
$ws.Range("H65").Value = 3875.3333
$ws.Range("I65").Value = 3104
$ws.Range("K65").Value = 15520
$ws.Range("M65").Value = -12400
$ws.Range("H100").Value = 1909.375
$ws.Range("I100").Value = 1640.8889
$ws.Range("J100").Value = 2254.5715
$ws.Range("K100").Value = 3281.7778
$ws.Range("L100").Value = 4509.143
$ws.Range("M100").Value = -2740.7778
$ws.Range("N100").Value = -5591.143
$ws.Range("H126").Value = 2082.5
$ws.Range("I126").Value = 1624.25
$ws.Range("J126").Value = 2999
$ws.Range("K126").Value = 4872.75
$ws.Range("L126").Value = 8997
$ws.Range("M126").Value = -2402.75
$ws.Range("N126").Value = -13937
$ws.Range("H132").Value = 3330.1875
$ws.Range("I132").Value = 2852.5386
$ws.Range("K132").Value = 8557.6158
$ws.Range("M132").Value = -6027.6158
